$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in C1
$ws.Range("C1").Value = "column3"

# Move the numeric value (and its number format) from B2 to C2
$numberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C2").Value = 123.123
$ws.Range("C2").NumberFormat = $numberFormat
$ws.Range("B2").Clear()

# Update selection to C1
$ws.Range("C1").Select()
